$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.112636208534241
$ws.Range("B1").Value = 1.041699409484863
$ws.Range("C1").Value = 5.165627956390381
$ws.Range("D1").Value = 1.62781822681427
$ws.Range("E1").Value = 0.9435285329818726
